$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.573880156280049
$ws.Range("C2").Value = 5.602865536938483
$ws.Range("D2").Value = 9.149089020520927
$ws.Range("F2").Value = 38.7198539006798
$ws.Range("G2").Value = 43.40943992380495
$ws.Range("H2").Value = 17.9367289136751
$ws.Range("I2").Value = 26.38136426604877
$ws.Range("J2").Value = 11.39450717807848
$ws.Range("K2").Value = 10.21435858003412
$ws.Range("M2").Value = 16.37561203550864
$ws.Range("B3").Value = 9.35209492403067
$ws.Range("C3").Value = 5.447729625047077
$ws.Range("D3").Value = 9.1205218787664
$ws.Range("F3").Value = 38.73683118239558
$ws.Range("G3").Value = 43.40280858259702
$ws.Range("H3").Value = 17.97652631681673
$ws.Range("I3").Value = 26.44666459872368
$ws.Range("J3").Value = 11.411592476549
$ws.Range("K3").Value = 10.07733013588103
$ws.Range("M3").Value = 16.32025643064994
$ws.Range("B4").Value = 9.21539884434087
$ws.Range("C4").Value = 5.35160740039935
$ws.Range("D4").Value = 9.104505924349846
$ws.Range("F4").Value = 38.75570600723574
$ws.Range("G4").Value = 43.40997534539664
$ws.Range("H4").Value = 18.00379710047492
$ws.Range("I4").Value = 26.49143600233822
$ws.Range("J4").Value = 11.42366484970966
$ws.Range("K4").Value = 9.994460492699623
$ws.Range("M4").Value = 16.28923429338353
$ws.Range("B5").Value = 9.159650751304003
$ws.Range("C5").Value = 5.312283677621309
$ws.Range("D5").Value = 9.098367411960114
$ws.Range("F5").Value = 38.76552101824889
$ws.Range("G5").Value = 43.41571916803468
$ws.Range("H5").Value = 18.01562234630585
$ws.Range("I5").Value = 26.51085462328827
$ws.Range("J5").Value = 11.42898228282968
$ws.Range("K5").Value = 9.961050054186343
$ws.Range("M5").Value = 16.27734792884798
$ws.Range("B6").Value = 9.150393809554519
$ws.Range("C6").Value = 5.305746723805426
$ws.Range("D6").Value = 9.097371691102895
$ws.Range("F6").Value = 38.76727897251204
$ws.Range("G6").Value = 43.41684333049674
$ws.Range("H6").Value = 18.01762891156089
$ws.Range("I6").Value = 26.51414990868195
$ws.Range("J6").Value = 11.42988926801095
$ws.Range("K6").Value = 9.955525188924922
$ws.Range("M6").Value = 16.27542009337536
$ws.Range("B7").Value = 9.214647064194988
$ws.Range("C7").Value = 5.351077598402024
$ws.Range("D7").Value = 9.104421560674067
$ws.Range("F7").Value = 38.75582978141654
$ws.Range("G7").Value = 43.41004138216721
$ws.Range("H7").Value = 18.00395369728646
$ws.Range("I7").Value = 26.4916931384095
$ws.Range("J7").Value = 11.42373495158363
$ws.Range("K7").Value = 9.994008397774449
$ws.Range("M7").Value = 16.28907091916386
$ws.Range("B8").Value = 9.497567709837671
$ws.Range("C8").Value = 5.549592891483698
$ws.Range("D8").Value = 9.138925671380123
$ws.Range("F8").Value = 38.72395302390976
$ws.Range("G8").Value = 43.4048198420973
$ws.Range("H8").Value = 17.94986223954174
$ws.Range("I8").Value = 26.4029079347635
$ws.Range("J8").Value = 11.40006990292106
$ws.Range("K8").Value = 10.16687317683588
$ws.Range("M8").Value = 16.35591565677894
$ws.Range("B9").Value = 10.04459017648775
$ws.Range("C9").Value = 5.929311161367052
$ws.Range("D9").Value = 9.218457077644461
$ws.Range("F9").Value = 38.72853057930801
$ws.Range("G9").Value = 43.48379249461903
$ws.Range("H9").Value = 17.86631407843202
$ws.Range("I9").Value = 26.265998843166
$ws.Range("J9").Value = 11.36621155832952
$ws.Range("K9").Value = 10.51408302654227
$ws.Range("M9").Value = 16.51008527837675
$ws.Range("B10").Value = 10.43710828600512
$ws.Range("C10").Value = 6.199144613771653
$ws.Range("D10").Value = 9.283799439800044
$ws.Range("F10").Value = 38.77278485172452
$ws.Range("G10").Value = 43.59613793942766
$ws.Range("H10").Value = 17.81870284171613
$ws.Range("I10").Value = 26.18820837631289
$ws.Range("J10").Value = 11.34898035189276
$ws.Range("K10").Value = 10.77174570847452
$ws.Range("M10").Value = 16.6367785521331
$ws.Range("B11").Value = 10.61274683983396
$ws.Range("C11").Value = 6.319296562676098
$ws.Range("D11").Value = 9.314949859705253
$ws.Range("F11").Value = 38.80177314143726
$ws.Range("G11").Value = 43.65898762745465
$ws.Range("H11").Value = 17.80004103216369
$ws.Range("I11").Value = 26.157792959669
$ws.Range("J11").Value = 11.34279928045811
$ws.Range("K11").Value = 10.88900971501156
$ws.Range("M11").Value = 16.6971768052423
$ws.Range("B12").Value = 10.6787690346125
$ws.Range("C12").Value = 6.364376035256083
$ws.Range("D12").Value = 9.326943825729646
$ws.Range("F12").Value = 38.81401989369622
$ws.Range("G12").Value = 43.68446759477232
$ws.Range("H12").Value = 17.79340572009519
$ws.Range("I12").Value = 26.14699209515747
$ws.Range("J12").Value = 11.34069673499098
$ws.Range("K12").Value = 10.93338042058788
$ws.Range("M12").Value = 16.72043170937325
$ws.Range("B13").Value = 10.66457272749437
$ws.Range("C13").Value = 6.354686708636619
$ws.Range("D13").Value = 9.324352020743934
$ws.Range("F13").Value = 38.81132595135239
$ws.Range("G13").Value = 43.67890546017578
$ws.Range("H13").Value = 17.79481555331845
$ws.Range("I13").Value = 26.14928634705794
$ws.Range("J13").Value = 11.34113897118921
$ws.Range("K13").Value = 10.92382659330377
$ws.Range("M13").Value = 16.71540652997497
$ws.Range("B14").Value = 10.61818868609836
$ws.Range("C14").Value = 6.323013933534656
$ws.Range("D14").Value = 9.315932684144999
$ws.Range("F14").Value = 38.80275527350195
$ws.Range("G14").Value = 43.66105026303111
$ws.Range("H14").Value = 17.79948648989383
$ws.Range("I14").Value = 26.15688999035032
$ws.Range("J14").Value = 11.34262153270404
$ws.Range("K14").Value = 10.8926610506983
$ws.Range("M14").Value = 16.69908239965512
$ws.Range("B15").Value = 10.58971158372029
$ws.Range("C15").Value = 6.303557522747758
$ws.Range("D15").Value = 9.310801165742955
$ws.Range("F15").Value = 38.79767067397648
$ws.Range("G15").Value = 43.65033195692097
$ws.Range("H15").Value = 17.80240378466466
$ws.Range("I15").Value = 26.16164084154904
$ws.Range("J15").Value = 11.34356064260191
$ws.Range("K15").Value = 10.87356551644572
$ws.Range("M15").Value = 16.68913289719018
$ws.Range("B16").Value = 10.42556517233286
$ws.Range("C16").Value = 6.191236115280144
$ws.Range("D16").Value = 9.281791781523841
$ws.Range("F16").Value = 38.77106829712752
$ws.Range("G16").Value = 43.59226618852127
$ws.Range("H16").Value = 17.81998278869134
$ws.Range("I16").Value = 26.19029630279954
$ws.Range("J16").Value = 11.34941763204315
$ws.Range("K16").Value = 10.76407965095295
$ws.Range("M16").Value = 16.63288582080639
$ws.Range("B17").Value = 10.32406851748412
$ws.Range("C17").Value = 6.121631378945621
$ws.Range("D17").Value = 9.264355727813308
$ws.Range("F17").Value = 38.75701449333936
$ws.Range("G17").Value = 43.55964704864836
$ws.Range("H17").Value = 17.83153491884167
$ws.Range("I17").Value = 26.20915026950222
$ws.Range("J17").Value = 11.35343506613738
$ws.Range("K17").Value = 10.69689662418419
$ws.Range("M17").Value = 16.5990782999351
$ws.Range("B18").Value = 10.26541926504374
$ws.Range("C18").Value = 6.081354648664897
$ws.Range("D18").Value = 9.254461619848518
$ws.Range("F18").Value = 38.74976507644181
$ws.Range("G18").Value = 43.54199134810737
$ws.Range("H18").Value = 17.83846143964768
$ws.Range("I18").Value = 26.22046244518182
$ws.Range("J18").Value = 11.3559018174775
$ws.Range("K18").Value = 10.6582620698223
$ws.Range("M18").Value = 16.57989423634438
$ws.Range("B19").Value = 10.24551724873338
$ws.Range("C19").Value = 6.067677562009377
$ws.Range("D19").Value = 9.251134981369532
$ws.Range("F19").Value = 38.74745388772268
$ws.Range("G19").Value = 43.53620360674622
$ws.Range("H19").Value = 17.8408550626921
$ws.Range("I19").Value = 26.22437285261067
$ws.Range("J19").Value = 11.35676382211674
$ws.Range("K19").Value = 10.64518363615382
$ws.Range("M19").Value = 16.57344411365453
$ws.Range("B20").Value = 10.33490157520756
$ws.Range("C20").Value = 6.129066307646127
$ws.Range("D20").Value = 9.266197940069729
$ws.Range("F20").Value = 38.75842425902137
$ws.Range("G20").Value = 43.56300500330759
$ws.Range("H20").Value = 17.83027598024103
$ws.Range("I20").Value = 26.20709479564064
$ws.Range("J20").Value = 11.35299125660584
$ws.Range("K20").Value = 10.70404792037188
$ws.Range("M20").Value = 16.60265024335565
$ws.Range("B21").Value = 10.63182657646212
$ws.Range("C21").Value = 6.332328725640397
$ws.Range("D21").Value = 9.318400331947396
$ws.Range("F21").Value = 38.80523827314266
$ws.Range("G21").Value = 43.66624924645002
$ws.Range("H21").Value = 17.7981028072939
$ws.Range("I21").Value = 26.15463714839525
$ws.Range("J21").Value = 11.34217960899073
$ws.Range("K21").Value = 10.90181639560715
$ws.Range("M21").Value = 16.70386690305619
$ws.Range("B22").Value = 10.82300978487987
$ws.Range("C22").Value = 6.462709099539271
$ws.Range("D22").Value = 9.35366811582125
$ws.Range("F22").Value = 38.84323096128822
$ws.Range("G22").Value = 43.74351313902044
$ws.Range("H22").Value = 17.77959116663465
$ws.Range("I22").Value = 26.12453182311913
$ws.Range("J22").Value = 11.33650119840843
$ws.Range("K22").Value = 11.03084901546644
$ws.Range("M22").Value = 16.77224559607034
$ws.Range("B23").Value = 10.72125602374229
$ws.Range("C23").Value = 6.393362172619807
$ws.Range("D23").Value = 9.33474216689117
$ws.Range("F23").Value = 38.82227835953695
$ws.Range("G23").Value = 43.70138370780462
$ws.Range("H23").Value = 17.78924084776841
$ws.Range("I23").Value = 26.14021666926172
$ws.Range("J23").Value = 11.33940500320282
$ws.Range("K23").Value = 10.96201565437601
$ws.Range("M23").Value = 16.73555159297197
$ws.Range("B24").Value = 10.3300048711131
$ws.Range("C24").Value = 6.125705782802672
$ws.Range("D24").Value = 9.265364670103375
$ws.Range("F24").Value = 38.75778431725684
$ws.Range("G24").Value = 43.5614834527266
$ws.Range("H24").Value = 17.83084425848093
$ws.Range("I24").Value = 26.20802260277993
$ws.Range("J24").Value = 11.35319141359953
$ws.Range("K24").Value = 10.70081484770972
$ws.Range("M24").Value = 16.60103458067777
$ws.Range("B25").Value = 9.897934728900291
$ws.Range("C25").Value = 5.82796987254564
$ws.Range("D25").Value = 9.19570325997246
$ws.Range("F25").Value = 38.72010654187955
$ws.Range("G25").Value = 43.45287583761079
$ws.Range("H25").Value = 17.88650026446518
$ws.Range("I25").Value = 26.29904123433396
$ws.Range("J25").Value = 11.37402808112355
$ws.Range("K25").Value = 10.41952685350155
$ws.Range("M25").Value = 16.46597087326617
